$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1395.3334
$ws.Range("I28").Value = 74.40000000000001
$ws.Range("J28").Value = 8000
$ws.Range("K28").Value = 74.40000000000001
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 410.6
$ws.Range("N28").Value = -8970

$ws.Range("H69").Value = 15759.467
$ws.Range("I69").Value = 6299
$ws.Range("J69").Value = 17214.924
$ws.Range("K69").Value = 18897
$ws.Range("L69").Value = 51644.772
$ws.Range("M69").Value = -18023
$ws.Range("N69").Value = -53392.772

$ws.Range("H72").Value = 15759.467
$ws.Range("I72").Value = 6299
$ws.Range("J72").Value = 17214.924
$ws.Range("K72").Value = 56691
$ws.Range("L72").Value = 154934.316
$ws.Range("M72").Value = -52323
$ws.Range("N72").Value = -163670.316

$ws.Range("H134").Value = 118330.664
$ws.Range("J134").Value = 118330.664
$ws.Range("L134").Value = 118330.664
$ws.Range("N134").Value = -128470.664

$ws.Range("H138").Value = 4190.8696
$ws.Range("I138").Value = 2405.24
$ws.Range("J138").Value = 5205.4316
$ws.Range("K138").Value = 7215.719999999999
$ws.Range("L138").Value = 15616.2948
$ws.Range("M138").Value = -2075.719999999999
$ws.Range("N138").Value = -25896.2948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1657.5
$ws.Range("I32").Value = 1509.644
$ws.Range("K32").Value = 1509.644
$ws.Range("M32").Value = -1222.644

$ws.Range("H102").Value = 3705035.2
$ws.Range("I102").Value = 4001309.2
$ws.Range("K102").Value = 4001309.2
$ws.Range("M102").Value = -3999687.2

$ws.Range("H132").Value = 3588158.2
$ws.Range("I132").Value = 2634777.5
$ws.Range("K132").Value = 7904332.5
$ws.Range("M132").Value = -7901802.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 10411.375
$ws.Range("J95").Value = 10411.375
$ws.Range("L95").Value = 10411.375
$ws.Range("N95").Value = -15903.375

$ws.Range("H99").Value = 2075.75
$ws.Range("I99").Value = 1951.9166
$ws.Range("J99").Value = 2447.25
$ws.Range("K99").Value = 1951.9166
$ws.Range("L99").Value = 2447.25
$ws.Range("M99").Value = -453.9166
$ws.Range("N99").Value = -5443.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7351
$ws.Range("I86").Value = 4921.2
$ws.Range("K86").Value = 4921.2
$ws.Range("M86").Value = -3798.2

$ws.Range("H89").Value = 7351
$ws.Range("I89").Value = 4921.2
$ws.Range("K89").Value = 24606
$ws.Range("M89").Value = -18990

$ws.Range("H92").Value = 38662.668
$ws.Range("J92").Value = 38662.668
$ws.Range("L92").Value = 38662.668
$ws.Range("N92").Value = -43654.668

$ws.Range("H132").Value = 35715636
$ws.Range("I132").Value = 37038404
$ws.Range("K132").Value = 111115212
$ws.Range("M132").Value = -111112682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 5690.5
$ws.Range("I82").Value = 5754
$ws.Range("J82").Value = 5500
$ws.Range("K82").Value = 17262
$ws.Range("L82").Value = 16500
$ws.Range("M82").Value = -16856
$ws.Range("N82").Value = -17312

$ws.Range("H85").Value = 5690.5
$ws.Range("I85").Value = 5754
$ws.Range("J85").Value = 5500
$ws.Range("K85").Value = 17262
$ws.Range("L85").Value = 16500
$ws.Range("M85").Value = -15858
$ws.Range("N85").Value = -19308

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H131").Value = 2128.5715
$ws.Range("I131").Value = 1812.4546
$ws.Range("J131").Value = 2476.3
$ws.Range("K131").Value = 5437.3638
$ws.Range("L131").Value = 7428.900000000001
$ws.Range("M131").Value = -397.3638000000001
$ws.Range("N131").Value = -17508.9

$ws.Range("H137").Value = 5266008
$ws.Range("J137").Value = 3272.3635
$ws.Range("L137").Value = 9817.0905
$ws.Range("N137").Value = -20017.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 12832.667
$ws.Range("J24").Value = 4500
$ws.Range("L24").Value = 4500
$ws.Range("N24").Value = -4846

$ws.Range("H53").Value = 22666.334
$ws.Range("I53").Value = 9000
$ws.Range("K53").Value = 9000
$ws.Range("M53").Value = -8369

$ws.Range("H70").Value = 6878.4
$ws.Range("I70").Value = 6848
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 6848
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -6578
$ws.Range("N70").Value = -7540

$ws.Range("H73").Value = 6878.4
$ws.Range("I73").Value = 6848
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 6848
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -5912
$ws.Range("N73").Value = -8872

$ws.Range("H132").Value = 2669362.8
$ws.Range("I132").Value = 3215006
$ws.Range("K132").Value = 9645018
$ws.Range("M132").Value = -9642488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1678.1
$ws.Range("I46").Value = 1684.7142
$ws.Range("K46").Value = 1684.7142
$ws.Range("M46").Value = -1496.7142

$ws.Range("H61").Value = 5991.0625
$ws.Range("I61").Value = 6057.1333
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 6057.1333
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -5855.1333
$ws.Range("N61").Value = -5404

$ws.Range("H68").Value = 1789421.2
$ws.Range("I68").Value = 2503759.5
$ws.Range("J68").Value = 3575.75
$ws.Range("K68").Value = 2503759.5
$ws.Range("L68").Value = 3575.75
$ws.Range("M68").Value = -2503010.5
$ws.Range("N68").Value = -5073.75

$ws.Range("H71").Value = 1789421.2
$ws.Range("I71").Value = 2503759.5
$ws.Range("J71").Value = 3575.75
$ws.Range("K71").Value = 12518797.5
$ws.Range("L71").Value = 17878.75
$ws.Range("M71").Value = -12515053.5
$ws.Range("N71").Value = -25366.75

$ws.Range("H113").Value = 5991.0625
$ws.Range("I113").Value = 6057.1333
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 6057.1333
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -3887.1333
$ws.Range("N113").Value = -9340

$ws.Range("H122").Value = 2998
$ws.Range("I122").Value = 2998
$ws.Range("K122").Value = 8994
$ws.Range("M122").Value = -6544

$ws.Range("H136").Value = 2135.8635
$ws.Range("I136").Value = 2103.6365
$ws.Range("K136").Value = 6310.9095
$ws.Range("M136").Value = -3760.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1567.2858
$ws.Range("I81").Value = 1567.2858
$ws.Range("K81").Value = 3134.5716
$ws.Range("M81").Value = -2073.5716

$ws.Range("H84").Value = 1567.2858
$ws.Range("I84").Value = 1567.2858
$ws.Range("K84").Value = 15672.858
$ws.Range("M84").Value = -10368.858

$ws.Range("H107").Value = 1428.7142
$ws.Range("I107").Value = 829.2727
$ws.Range("K107").Value = 2487.8181
$ws.Range("M107").Value = -567.8181

$ws.Range("H113").Value = 1498.1428
$ws.Range("I113").Value = 1477.4
$ws.Range("J113").Value = 1550
$ws.Range("K113").Value = 4432.200000000001
$ws.Range("L113").Value = 4650
$ws.Range("M113").Value = -2262.200000000001
$ws.Range("N113").Value = -8990

$ws.Range("H132").Value = 8476253
$ws.Range("I132").Value = 10418112
$ws.Range("J132").Value = 2685.5454
$ws.Range("K132").Value = 31254336
$ws.Range("L132").Value = 8056.6362
$ws.Range("M132").Value = -31251806
$ws.Range("N132").Value = -13116.6362

